# Scheduled data refresh: update market-price snapshot columns (H:N)
# across the per-job worksheets (ALC/ARM/BSM/CRP/GSM/LTW/WVR), row by row.
$wb = $excel.ActiveWorkbook

# ALC!98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1382.7894
$ws.Range("I98").Value = 784.1429000000001
$ws.Range("J98").Value = 3059
$ws.Range("K98").Value = 784.1429000000001
$ws.Range("L98").Value = 3059
$ws.Range("M98").Value = 713.8570999999999
$ws.Range("N98").Value = -6055

# ALC!122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1382.7894
$ws.Range("I122").Value = 784.1429000000001
$ws.Range("J122").Value = 3059
$ws.Range("K122").Value = 2352.4287
$ws.Range("L122").Value = 9177
$ws.Range("M122").Value = 97.57129999999961
$ws.Range("N122").Value = -14077

# ALC!132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2207.494
$ws.Range("I132").Value = 1615.6307
$ws.Range("J132").Value = 4344.778
$ws.Range("K132").Value = 4846.8921
$ws.Range("L132").Value = 13034.334
$ws.Range("M132").Value = -2316.8921
$ws.Range("N132").Value = -18094.334

# ARM!32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4129.2627
$ws.Range("I32").Value = 3143.1572
$ws.Range("J32").Value = 12905.6
$ws.Range("K32").Value = 3143.1572
$ws.Range("L32").Value = 12905.6
$ws.Range("M32").Value = -2856.1572
$ws.Range("N32").Value = -13479.6

# ARM!61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1019.0492
$ws.Range("I61").Value = 791.87756
$ws.Range("J61").Value = 1946.6666
$ws.Range("K61").Value = 791.87756
$ws.Range("L61").Value = 1946.6666
$ws.Range("M61").Value = -579.87756
$ws.Range("N61").Value = -2370.6666

# ARM!88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2619.6
$ws.Range("I88").Value = 2231
$ws.Range("J88").Value = 3040.5833
$ws.Range("K88").Value = 2231
$ws.Range("L88").Value = 3040.5833
$ws.Range("M88").Value = -1825
$ws.Range("N88").Value = -3852.5833

# ARM!91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2619.6
$ws.Range("I91").Value = 2231
$ws.Range("J91").Value = 3040.5833
$ws.Range("K91").Value = 2231
$ws.Range("L91").Value = 3040.5833
$ws.Range("M91").Value = -827
$ws.Range("N91").Value = -5848.5833

# ARM!92
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 32399.6
$ws.Range("J92").Value = 32399.6
$ws.Range("L92").Value = 32399.6
$ws.Range("N92").Value = -37391.6

# ARM!132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5209.98
$ws.Range("I132").Value = 3473.675
$ws.Range("K132").Value = 10421.025
$ws.Range("M132").Value = -7891.025000000001

# ARM!136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1019.0492
$ws.Range("I136").Value = 791.87756
$ws.Range("J136").Value = 1946.6666
$ws.Range("K136").Value = 2375.63268
$ws.Range("L136").Value = 5839.9998
$ws.Range("M136").Value = 174.3673199999998
$ws.Range("N136").Value = -10939.9998

# BSM!86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2006.0312
$ws.Range("I86").Value = 1800.9375
$ws.Range("J86").Value = 2211.125
$ws.Range("K86").Value = 1800.9375
$ws.Range("L86").Value = 2211.125
$ws.Range("M86").Value = -677.9375
$ws.Range("N86").Value = -4457.125

# BSM!89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2006.0312
$ws.Range("I89").Value = 1800.9375
$ws.Range("J89").Value = 2211.125
$ws.Range("K89").Value = 9004.6875
$ws.Range("L89").Value = 11055.625
$ws.Range("M89").Value = -3388.6875
$ws.Range("N89").Value = -22287.625

# BSM!134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1179.3636
$ws.Range("I134").Value = 856.5714
$ws.Range("J134").Value = 1744.25
$ws.Range("K134").Value = 2569.7142
$ws.Range("L134").Value = 5232.75
$ws.Range("M134").Value = -34.71420000000035
$ws.Range("N134").Value = -10302.75

# CRP!17
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

# CRP!31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13187554
$ws.Range("I31").Value = 19231838
$ws.Range("J31").Value = 91607.086
$ws.Range("K31").Value = 19231838
$ws.Range("L31").Value = 91607.086
$ws.Range("M31").Value = -19231543
$ws.Range("N31").Value = -92197.086

# CRP!34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 13187554
$ws.Range("I34").Value = 19231838
$ws.Range("J34").Value = 91607.086
$ws.Range("K34").Value = 19231838
$ws.Range("L34").Value = 91607.086
$ws.Range("M34").Value = -19231636
$ws.Range("N34").Value = -92011.086

# CRP!58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1171.9375
$ws.Range("I58").Value = 845.5161000000001
$ws.Range("J58").Value = 1767.1765
$ws.Range("K58").Value = 845.5161000000001
$ws.Range("L58").Value = 1767.1765
$ws.Range("M58").Value = -642.5161000000001
$ws.Range("N58").Value = -2173.1765

# CRP!62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4094.5676
$ws.Range("I62").Value = 4580.769
$ws.Range("J62").Value = 2945.3635
$ws.Range("K62").Value = 4580.769
$ws.Range("L62").Value = 2945.3635
$ws.Range("M62").Value = -3956.769
$ws.Range("N62").Value = -4193.363499999999

# CRP!65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4094.5676
$ws.Range("I65").Value = 4580.769
$ws.Range("J65").Value = 2945.3635
$ws.Range("K65").Value = 22903.845
$ws.Range("L65").Value = 14726.8175
$ws.Range("M65").Value = -19783.845
$ws.Range("N65").Value = -20966.8175

# CRP!122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1211.15
$ws.Range("I122").Value = 1141.5333
$ws.Range("J122").Value = 1420
$ws.Range("K122").Value = 3424.5999
$ws.Range("L122").Value = 4260
$ws.Range("M122").Value = -974.5999000000002
$ws.Range("N122").Value = -9160

# CRP!136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1171.9375
$ws.Range("I136").Value = 845.5161000000001
$ws.Range("J136").Value = 1767.1765
$ws.Range("K136").Value = 2536.5483
$ws.Range("L136").Value = 5301.529500000001
$ws.Range("M136").Value = 13.45169999999962
$ws.Range("N136").Value = -10401.5295

# GSM!102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3327.8
$ws.Range("I102").Value = 3284.35
$ws.Range("J102").Value = 3501.6
$ws.Range("K102").Value = 3284.35
$ws.Range("L102").Value = 3501.6
$ws.Range("M102").Value = -1662.35
$ws.Range("N102").Value = -6745.6

# LTW!13
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 1006
$ws.Range("I13").Value = 1006
$ws.Range("K13").Value = 1006
$ws.Range("M13").Value = -866

# LTW!68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1987.0435
$ws.Range("I68").Value = 2007.2858
$ws.Range("J68").Value = 1955.5555
$ws.Range("K68").Value = 2007.2858
$ws.Range("L68").Value = 1955.5555
$ws.Range("M68").Value = -1258.2858
$ws.Range("N68").Value = -3453.5555

# LTW!71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1987.0435
$ws.Range("I71").Value = 2007.2858
$ws.Range("J71").Value = 1955.5555
$ws.Range("K71").Value = 10036.429
$ws.Range("L71").Value = 9777.7775
$ws.Range("M71").Value = -6292.429
$ws.Range("N71").Value = -17265.7775

# WVR!122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4763376.5
$ws.Range("I122").Value = 6667971.5
$ws.Range("J122").Value = 1888.75
$ws.Range("K122").Value = 20003914.5
$ws.Range("L122").Value = 5666.25
$ws.Range("M122").Value = -20001464.5
$ws.Range("N122").Value = -10566.25

# WVR!132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8336213
$ws.Range("I132").Value = 11908046
$ws.Range("J132").Value = 1936.1111
$ws.Range("K132").Value = 35724138
$ws.Range("L132").Value = 5808.3333
$ws.Range("M132").Value = -35721608
$ws.Range("N132").Value = -10868.3333
